$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 500
$ws.Range("I9").Value = 500
$ws.Range("K9").Value = 500
$ws.Range("M9").Value = -331

$ws.Range("H12").Value = 333.33334
$ws.Range("I12").Value = 250
$ws.Range("J12").Value = 500
$ws.Range("K12").Value = 250
$ws.Range("L12").Value = 500
$ws.Range("M12").Value = -80
$ws.Range("N12").Value = -840

$ws.Range("H21").Value = 35021
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 35021
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 35021
$ws.Range("N21").Value = -35957
$ws.Range("M21").ClearContents()

$ws.Range("H23").Value = 35021
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 35021
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 35021
$ws.Range("N23").Value = -35489
$ws.Range("M23").ClearContents()

$ws.Range("H29").Value = 180
$ws.Range("I29").Value = 180
$ws.Range("K29").Value = 540
$ws.Range("M29").Value = -259

$ws.Range("H39").Value = 237.75
$ws.Range("I39").Value = 237.75
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 713.25
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -417.25
$ws.Range("N39").ClearContents()

$ws.Range("H41").Value = 632.2222
$ws.Range("I41").Value = 694
$ws.Range("K41").Value = 694
$ws.Range("M41").Value = -254

$ws.Range("H88").Value = 2210.6667
$ws.Range("J88").Value = 2210.6667
$ws.Range("L88").Value = 2210.6667
$ws.Range("N88").Value = -3022.6667

$ws.Range("H91").Value = 2210.6667
$ws.Range("J91").Value = 2210.6667
$ws.Range("L91").Value = 2210.6667
$ws.Range("N91").Value = -5018.6667

$ws.Range("H97").Value = 1998
$ws.Range("J97").Value = 1998
$ws.Range("L97").Value = 5994
$ws.Range("N97").Value = -6986

$ws.Range("H98").Value = 2140.7144
$ws.Range("I98").Value = 2140.7144
$ws.Range("K98").Value = 2140.7144
$ws.Range("M98").Value = -642.7143999999998

$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()

$ws.Range("H112").Value = 1415
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H116").Value = 4071.2856
$ws.Range("I116").Value = 3999.8333
$ws.Range("K116").Value = 3999.8333
$ws.Range("M116").Value = -557.8332999999998

$ws.Range("H122").Value = 2140.7144
$ws.Range("I122").Value = 2140.7144
$ws.Range("K122").Value = 6422.1432
$ws.Range("M122").Value = -3972.1432

$ws.Range("H132").Value = 3430.4167
$ws.Range("I132").Value = 2836.6
$ws.Range("J132").Value = 6399.5
$ws.Range("K132").Value = 8509.799999999999
$ws.Range("L132").Value = 19198.5
$ws.Range("M132").Value = -5979.799999999999
$ws.Range("N132").Value = -24258.5

$ws.Range("H135").Value = 2698
$ws.Range("I135").Value = 2698
$ws.Range("K135").Value = 24282
$ws.Range("M135").Value = -21747

$ws.Range("H137").Value = 2438.3572
$ws.Range("I137").Value = 1648.1428
$ws.Range("K137").Value = 4944.428400000001
$ws.Range("M137").Value = -2394.428400000001

$ws.Range("H138").Value = 4499.375
$ws.Range("J138").Value = 5400
$ws.Range("L138").Value = 16200
$ws.Range("N138").Value = -26480

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29008.375
$ws.Range("I32").Value = 29008.375
$ws.Range("K32").Value = 29008.375
$ws.Range("M32").Value = -28721.375

$ws.Range("H45").Value = 2039.6
$ws.Range("I45").Value = 1499.6666
$ws.Range("K45").Value = 1499.6666
$ws.Range("M45").Value = -1122.6666

$ws.Range("H88").Value = 1660.25
$ws.Range("J88").Value = 1837.6
$ws.Range("L88").Value = 1837.6
$ws.Range("N88").Value = -2649.6

$ws.Range("H91").Value = 1660.25
$ws.Range("J91").Value = 1837.6
$ws.Range("L91").Value = 1837.6
$ws.Range("N91").Value = -4645.6

$ws.Range("H102").Value = 2483.8333
$ws.Range("I102").Value = 2483.8333
$ws.Range("K102").Value = 2483.8333
$ws.Range("M102").Value = -861.8332999999998

$ws.Range("H110").Value = 1284.8334
$ws.Range("I110").Value = 1284.8334
$ws.Range("K110").Value = 1284.8334
$ws.Range("M110").Value = 760.1666

$ws.Range("H132").Value = 2949.4375
$ws.Range("I132").Value = 2254.889
$ws.Range("K132").Value = 6764.667
$ws.Range("M132").Value = -4234.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1763.375
$ws.Range("I99").Value = 1758.1428
$ws.Range("K99").Value = 1758.1428
$ws.Range("M99").Value = -260.1428000000001

$ws.Range("H134").Value = 3685.16
$ws.Range("I134").Value = 3596.818
$ws.Range("J134").Value = 4333
$ws.Range("K134").Value = 10790.454
$ws.Range("L134").Value = 12999
$ws.Range("M134").Value = -8255.454000000002
$ws.Range("N134").Value = -18069

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4149.5625
$ws.Range("J31").Value = 13198
$ws.Range("L31").Value = 13198
$ws.Range("N31").Value = -13788

$ws.Range("H34").Value = 4149.5625
$ws.Range("J34").Value = 13198
$ws.Range("L34").Value = 13198
$ws.Range("N34").Value = -13602

$ws.Range("H58").Value = 2926
$ws.Range("I58").Value = 3047.8572
$ws.Range("J58").Value = 2499.5
$ws.Range("K58").Value = 3047.8572
$ws.Range("L58").Value = 2499.5
$ws.Range("M58").Value = -2844.8572
$ws.Range("N58").Value = -2905.5

$ws.Range("H99").Value = 6168.385
$ws.Range("J99").Value = 7498.6665
$ws.Range("L99").Value = 7498.6665
$ws.Range("N99").Value = -10494.6665

$ws.Range("H107").Value = 1789.16
$ws.Range("I107").Value = 1554.5
$ws.Range("K107").Value = 1554.5
$ws.Range("M107").Value = 365.5

$ws.Range("H126").Value = 6168.385
$ws.Range("J126").Value = 7498.6665
$ws.Range("L126").Value = 22495.9995
$ws.Range("N126").Value = -27435.9995

$ws.Range("H134").Value = 2143.3
$ws.Range("I134").Value = 1815.7059
$ws.Range("K134").Value = 5447.1177
$ws.Range("M134").Value = -2912.1177

$ws.Range("H136").Value = 2926
$ws.Range("I136").Value = 3047.8572
$ws.Range("J136").Value = 2499.5
$ws.Range("K136").Value = 9143.571599999999
$ws.Range("L136").Value = 7498.5
$ws.Range("M136").Value = -6593.571599999999
$ws.Range("N136").Value = -12598.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 66821784
$ws.Range("J4").Value = 26733334
$ws.Range("L4").Value = 80200002
$ws.Range("N4").Value = -80200226

$ws.Range("H8").Value = 450
$ws.Range("I8").Value = 450
$ws.Range("K8").Value = 1350
$ws.Range("M8").Value = -1211

$ws.Range("H122").Value = 201999.4
$ws.Range("I122").Value = 1500
$ws.Range("J122").Value = 252124.25
$ws.Range("K122").Value = 13500
$ws.Range("L122").Value = 2269118.25
$ws.Range("M122").Value = -11050
$ws.Range("N122").Value = -2274018.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2020.6
$ws.Range("I102").Value = 1774
$ws.Range("K102").Value = 1774
$ws.Range("M102").Value = -152

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1499.5
$ws.Range("J40").Value = 999
$ws.Range("L40").Value = 999
$ws.Range("N40").Value = -1271

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3609
$ws.Range("J132").Value = 4000.25
$ws.Range("L132").Value = 12000.75
$ws.Range("N132").Value = -17060.75
